# Added more RAD Test Cases and Data for MRF.
# Fills in Result/Date columns (A/B) for additional executed test rows
# across the Extension, NewTaxReturn, Personal_IND, Personal_JNT and
# Personal_EL sheets, and refreshes the latest run timestamp on the
# Estimated sheet.

$wb = $excel.ActiveWorkbook

# --- Estimated: update the most recent run's timestamp (row 2) ---
$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("B2").Value = "Tue Oct 03 12:04:00 EDT 2023"

# --- Extension: rows 2-7 ---
$ws = $wb.Worksheets.Item("Extension")
$dates = @(
    "Tue Oct 03 13:05:13 EDT 2023",
    "Tue Oct 03 13:05:35 EDT 2023",
    "Tue Oct 03 13:05:55 EDT 2023",
    "Tue Oct 03 13:06:14 EDT 2023",
    "Tue Oct 03 13:06:33 EDT 2023",
    "Tue Oct 03 13:06:52 EDT 2023"
)
$row = 2
foreach ($d in $dates) {
    $ws.Cells.Item($row, 1).Value = "Pass"
    $ws.Cells.Item($row, 2).Value = $d
    $row = $row + 1
}

# --- NewTaxReturn: rows 2-16 ---
$ws = $wb.Worksheets.Item("NewTaxReturn")
$dates = @(
    "Tue Oct 03 13:10:07 EDT 2023",
    "Tue Oct 03 13:10:28 EDT 2023",
    "Tue Oct 03 13:10:47 EDT 2023",
    "Tue Oct 03 13:11:07 EDT 2023",
    "Tue Oct 03 13:11:26 EDT 2023",
    "Tue Oct 03 13:11:45 EDT 2023",
    "Tue Oct 03 13:12:04 EDT 2023",
    "Tue Oct 03 13:12:24 EDT 2023",
    "Tue Oct 03 13:12:43 EDT 2023",
    "Tue Oct 03 13:13:01 EDT 2023",
    "Tue Oct 03 13:13:20 EDT 2023",
    "Tue Oct 03 13:13:39 EDT 2023",
    "Tue Oct 03 13:13:58 EDT 2023",
    "Tue Oct 03 13:14:17 EDT 2023",
    "Tue Oct 03 13:14:36 EDT 2023"
)
$row = 2
foreach ($d in $dates) {
    $ws.Cells.Item($row, 1).Value = "Pass"
    $ws.Cells.Item($row, 2).Value = $d
    $row = $row + 1
}

# --- Personal_EL: row 2 ---
$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Cells.Item(2, 1).Value = "Pass"
$ws.Cells.Item(2, 2).Value = "Tue Oct 03 15:28:36 EDT 2023"

# --- Personal_IND: rows 2-6 ---
$ws = $wb.Worksheets.Item("Personal_IND")
$dates = @(
    "Tue Oct 03 15:34:43 EDT 2023",
    "Tue Oct 03 15:35:04 EDT 2023",
    "Tue Oct 03 15:35:24 EDT 2023",
    "Tue Oct 03 15:35:44 EDT 2023",
    "Tue Oct 03 15:36:03 EDT 2023"
)
$row = 2
foreach ($d in $dates) {
    $ws.Cells.Item($row, 1).Value = "Pass"
    $ws.Cells.Item($row, 2).Value = $d
    $row = $row + 1
}

# --- Personal_JNT: rows 2-6 ---
$ws = $wb.Worksheets.Item("Personal_JNT")
$dates = @(
    "Tue Oct 03 15:45:13 EDT 2023",
    "Tue Oct 03 15:45:43 EDT 2023",
    "Tue Oct 03 15:46:10 EDT 2023",
    "Tue Oct 03 15:46:38 EDT 2023",
    "Tue Oct 03 15:47:05 EDT 2023"
)
$row = 2
foreach ($d in $dates) {
    $ws.Cells.Item($row, 1).Value = "Pass"
    $ws.Cells.Item($row, 2).Value = $d
    $row = $row + 1
}
